# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values on the zh-cn and de-de
# sheets for the d0e5e3d8... entry (row 2) with fresh report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 14:58:27"
$wsZhCn.Range("H2").Value = "2016-03-20 14:58:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 14:58:31"
$wsDeDe.Range("H2").Value = "2016-03-20 14:58:59"
